{"js": "// Office.js (Word JavaScript API) equivalent of the OOXML diff:\n// - The title paragraph's date text changes from \"2023-05-29 Monday\" to \"2023-05-30 Tuesday\".\n// - Each of the 100 arithmetic-answer cells in the single table is replaced with a new\n//   expression/result, preserving cell position (row-major) and all existing run formatting\n//   (font, size, etc.), since insertText(..., Replace) only swaps the text of the paragraph,\n//   leaving the paragraph/run properties untouched.\n//\n// The body contains exactly 101 paragraphs in document order: paragraph 0 is the title/date\n// line above the table, and paragraphs 1..100 are the single paragraph inside each of the\n// 100 table cells (20 rows x 5 columns), enumerated row-major - matching this list exactly.\nconst newValues = [\n  \"2023-05-30 Tuesday\", \"70+19=89\", \"24+74=98\", \"4+47=51\", \"32+22=54\", \"51-0=51\", \"83-32=51\",\n  \"13+42=55\", \"26+20=46\", \"85-28=57\", \"10-8=2\", \"47-23=24\", \"49+50=99\", \"22+69=91\", \"60+23=83\",\n  \"99-70=29\", \"0+34=34\", \"58-41=17\", \"15+20=35\", \"26+2=28\", \"62-28=34\", \"64-38=26\", \"54+2=56\",\n  \"39-34=5\", \"63+23=86\", \"96-3=93\", \"19-13=6\", \"71+20=91\", \"25+55=80\", \"68-17=51\", \"48-26=22\",\n  \"44+54=98\", \"42+10=52\", \"37+35=72\", \"73-68=5\", \"74-52=22\", \"14+80=94\", \"30-13=17\", \"6+70=76\",\n  \"44-37=7\", \"19+36=55\", \"9+29=38\", \"93-43=50\", \"22+41=63\", \"65-25=40\", \"58-39=19\", \"12+43=55\",\n  \"99-69=30\", \"26-3=23\", \"83-50=33\", \"35-16=19\", \"78-35=43\", \"4+41=45\", \"25+11=36\", \"58+21=79\",\n  \"72+18=90\", \"84-3=81\", \"11+38=49\", \"86-72=14\", \"65-31=34\", \"75-3=72\", \"96-76=20\", \"22+76=98\",\n  \"26+39=65\", \"12+5=17\", \"68+15=83\", \"34-17=17\", \"80-77=3\", \"0+92=92\", \"94-34=60\", \"3+18=21\",\n  \"76-58=18\", \"17+13=30\", \"51-50=1\", \"9+10=19\", \"58+9=67\", \"20+63=83\", \"34+5=39\", \"1-0=1\",\n  \"6+75=81\", \"64-47=17\", \"44-35=9\", \"79-12=67\", \"97-85=12\", \"88-82=6\", \"4+48=52\", \"20-16=4\",\n  \"37+51=88\", \"5+71=76\", \"94-43=51\", \"56-14=42\", \"25+70=95\", \"9+85=94\", \"4+24=28\", \"91-47=44\",\n  \"92-91=1\", \"69-11=58\", \"47-4=43\", \"42+20=62\", \"42+25=67\", \"34+22=56\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: \" + paragraphs.items.length + \" (expected \" + newValues.length + \")\"\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].insertText(newValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop equivalent of the OOXML diff:\n# - The title paragraph's date text changes from \"2023-05-29 Monday\" to \"2023-05-30 Tuesday\".\n# - Each of the 100 arithmetic-answer cells in the single table is replaced with a new\n#   expression/result, preserving cell position (row-major) and all existing run formatting\n#   (font, size, etc.), since assigning to Range.Text only swaps the text inside the existing\n#   run(s)/paragraph mark, leaving paragraph/run properties untouched.\n\n$d = $word.ActiveDocument\n\n# The title/date line is the first paragraph in the document body, above the table.\n$d.Paragraphs.Item(1).Range.Text = '2023-05-30 Tuesday'\n\n# The 100 answer values below are in row-major order (row 1 col 1..5, row 2 col 1..5, ...),\n# matching the single table's 20 rows x 5 columns exactly.\n$newValues = @(\n    '70+19=89', '24+74=98', '4+47=51', '32+22=54', '51-0=51', '83-32=51', '13+42=55', '26+20=46', \n    '85-28=57', '10-8=2', '47-23=24', '49+50=99', '22+69=91', '60+23=83', '99-70=29', '0+34=34', \n    '58-41=17', '15+20=35', '26+2=28', '62-28=34', '64-38=26', '54+2=56', '39-34=5', '63+23=86', \n    '96-3=93', '19-13=6', '71+20=91', '25+55=80', '68-17=51', '48-26=22', '44+54=98', '42+10=52', \n    '37+35=72', '73-68=5', '74-52=22', '14+80=94', '30-13=17', '6+70=76', '44-37=7', '19+36=55', \n    '9+29=38', '93-43=50', '22+41=63', '65-25=40', '58-39=19', '12+43=55', '99-69=30', '26-3=23', \n    '83-50=33', '35-16=19', '78-35=43', '4+41=45', '25+11=36', '58+21=79', '72+18=90', '84-3=81', \n    '11+38=49', '86-72=14', '65-31=34', '75-3=72', '96-76=20', '22+76=98', '26+39=65', '12+5=17', \n    '68+15=83', '34-17=17', '80-77=3', '0+92=92', '94-34=60', '3+18=21', '76-58=18', '17+13=30', \n    '51-50=1', '9+10=19', '58+9=67', '20+63=83', '34+5=39', '1-0=1', '6+75=81', '64-47=17', \n    '44-35=9', '79-12=67', '97-85=12', '88-82=6', '4+48=52', '20-16=4', '37+51=88', '5+71=76', \n    '94-43=51', '56-14=42', '25+70=95', '9+85=94', '4+24=28', '91-47=44', '92-91=1', '69-11=58', \n    '47-4=43', '42+20=62', '42+25=67', '34+22=56'\n)\n\n$t = $d.Tables.Item(1)\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
